$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 10. This shifts the old row 10 -> 11 and old row 11 -> 12,
# carrying all their existing values/styles along automatically.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record's data.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44776
$ws.Range("D10").NumberFormat = $ws.Range("D11").NumberFormat
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100107
$ws.Range("H10").Value = "Otros"
$ws.Range("I10").Value = 100107002
$ws.Range("J10").Value = "Chirimoya"
$ws.Range("K10").Value = "Cultivar IV Región"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 29000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 29500
$ws.Range("Q10").Value = "$/caja 10 kilos"
$ws.Range("R10").Value = "Región de Coquimbo"
$ws.Range("S10").Value = 2950
$ws.Range("T10").Value = 10
